$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.525.62"
$ws.Range("E2").Value = "  -0.01%  "

$ws.Range("D3").Value = "3.102.54"
$ws.Range("E3").Value = "  +2.92%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "385.54"
$ws.Range("E5").Value = "  +1.96%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.35"
$ws.Range("E6").Value = "  +0.41%  "

$ws.Range("E7").Value = "  -0.94%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  -1.78%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.09"
$ws.Range("E10").Value = "  +1.47%  "

$ws.Range("E11").Value = "  +0.07%  "

$ws.Range("E12").Value = "  -0.41%  "

$ws.Range("D13").Value = "3.589.79"
$ws.Range("E13").Value = "  +2.64%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.60"
$ws.Range("E14").Value = "  +0.78%  "

$ws.Range("E15").Value = "  +1.45%  "

$ws.Range("D16").Value = "3.104.45"
$ws.Range("E16").Value = "  +2.57%  "

$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.995"
$ws.Range("E17").Value = "  +1.70%  "

$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.00"
$ws.Range("E18").Value = "  +4.58%  "

$ws.Range("D19").Value = "51.555.98"
$ws.Range("E19").Value = "  +0.04%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.28"
$ws.Range("E20").Value = "  +8.44%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.37"
$ws.Range("E21").Value = "  -0.82%  "

$ws.Range("D22").Value = "0.0₃0964"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.88"
$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "266.19"
$ws.Range("E24").Value = "  -0.37%  "

$ws.Range("E25").Value = "  +0.66%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.09"
$ws.Range("E26").Value = "  -1.35%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.05"
$ws.Range("E27").Value = "  +3.24%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.26"
$ws.Range("E28").Value = "  -3.36%  "

$ws.Range("E29").Value = "  +0.04%  "

$ws.Range("E30").Value = "  -2.77%  "

$ws.Range("E31").Value = "  -1.58%  "

$ws.Range("E32").Value = "  +0.64%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0475"
$ws.Range("E33").Value = "  +4.29%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "35.21"
$ws.Range("E34").Value = "  +3.22%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.06"
$ws.Range("E35").Value = "  +0.67%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "50.29"
$ws.Range("E36").Value = "  -0.77%  "

$ws.Range("E37").Value = "  -0.21%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.36"
$ws.Range("E38").Value = "  +2.07%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.293"
$ws.Range("E39").Value = "  +2.35%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.89"
$ws.Range("E40").Value = "  +1.70%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "128.89"
$ws.Range("E41").Value = "  +4.81%  "

$ws.Range("E42").Value = "  -0.30%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.55"
$ws.Range("E43").Value = "  -3.72%  "

$ws.Range("E44").Value = "  -3.37%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "22.55"
$ws.Range("E45").Value = "  +5.07%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.67"
$ws.Range("E46").Value = "  -1.73%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.49"
$ws.Range("E47").Value = "  +4.73%  "

$ws.Range("E48").Value = "  +0.35%  "

$ws.Range("D49").Value = "2.065.68"
$ws.Range("E49").Value = "  +1.72%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0332"
$ws.Range("E50").Value = "  +4.06%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.917"
$ws.Range("E51").Value = "  +15.92%  "
